# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) and DialogAct (col J)
# values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 7;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 8;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 11;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 27;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 34;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 37;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 54;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 56;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 60;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 62;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 114; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 126; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 129; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 138; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 150; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 151; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 176; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 183; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 185; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 188; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 193; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 194; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 196; Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 197; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
